# Add a new project row (row 28) to the "2024-2" sheet, copying the
# formatting of the last existing row (27) and filling in the new
# E/P MALAGA 1 data, then restore the view/selection state that was
# saved with the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024-2")

# Copy formatting from the row above so the new row matches the
# existing table styling (borders, fonts, number formats, etc.)
$ws.Range("A27:G27").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values
$ws.Range("A28").Value = "2024-2"
$ws.Range("B28").Value = "E/P MALAGA 1"
$ws.Range("C28").Value = "Embarcación Pesquera"
$ws.Range("D28").Value = "A.S/0031"
$ws.Range("E28").Value = "A.S/0031-224"
$ws.Range("F28").Value = 45594
$ws.Range("G28").Value = 45600

# Restore the sheet view state saved with the workbook
$ws.Activate()
$ws.Range("L22").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 2
